$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.410.74'
$ws.Range("E2").Value = '  +4.14%  '

$ws.Range("D3").Value = '3.132.21'
$ws.Range("E3").Value = '  +2.44%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.49%  '

$ws.Range("E7").Value = '  +3.81%  '

$ws.Range("E8").Value = '  +20.51%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '3.130.07'
$ws.Range("E10").Value = '  +2.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.735'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +23.26%  '

$ws.Range("E12").Value = '  +6.44%  '

$ws.Range("E13").Value = '  +8.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +12.49%  '

$ws.Range("E15").Value = '  +3.91%  '

$ws.Range("D16").Value = '91.233.58'
$ws.Range("E16").Value = '  +4.02%  '

$ws.Range("D17").Value = '3.712.48'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("D18").Value = '3.133.36'
$ws.Range("E18").Value = '  +2.76%  '

$ws.Range("E19").Value = '  +19.23%  '

$ws.Range("E20").Value = '  +15.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.80%  '

$ws.Range("E23").Value = '  +9.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '86.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.64%  '

$ws.Range("D28").Value = '3.294.41'
$ws.Range("E28").Value = '  +2.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  +4.50%  '

$ws.Range("E31").Value = '  +12.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '532.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.16%  '

$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.63%  '

$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.900'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -17.07%  '

$ws.Range("E35").Value = '  +10.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.144'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.85'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.60%  '

$ws.Range("E38").Value = '  +6.23%  '

$ws.Range("E39").Value = '  +4.80%  '

$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("E42").Value = '  +16.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0791'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +20.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.383'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.59%  '

$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.72%  '

$ws.Range("E50").Value = '  +9.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000260'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +24.54%  '
